$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 36373.28826856845
$ws.Range("C2").Value = 32828.32581981639
$ws.Range("D2").Value = 33111.22475746145
$ws.Range("E2").Value = 33427.38324280769
$ws.Range("F2").Value = 33774.83582826793
$ws.Range("G2").Value = 33963.38123296577
$ws.Range("H2").Value = 33886.06427502663
$ws.Range("I2").Value = 33818.51708178314
$ws.Range("J2").Value = 32766.82248902799
$ws.Range("K2").Value = 28689.97216651863
$ws.Range("B3").Value = 28644.23911667487
$ws.Range("C3").Value = 25437.85421204317
$ws.Range("D3").Value = 25372.37373150306
$ws.Range("E3").Value = 25320.62365315176
$ws.Range("F3").Value = 25281.80043190278
$ws.Range("G3").Value = 25130.52850649365
$ws.Range("H3").Value = 24796.61272205292
$ws.Range("I3").Value = 24463.25406783813
$ws.Range("J3").Value = 23146.60140917578
$ws.Range("K3").Value = 23125.19029400865
$ws.Range("B4").Value = 28159.8542017977
$ws.Range("C4").Value = 28159.8542017977
$ws.Range("D4").Value = 28159.8542017977
$ws.Range("E4").Value = 28159.8542017977
$ws.Range("F4").Value = 28159.8542017977
$ws.Range("G4").Value = 28159.8542017977
$ws.Range("H4").Value = 28159.8542017977
$ws.Range("I4").Value = 28159.8542017977
$ws.Range("J4").Value = 28159.8542017977
$ws.Range("K4").Value = 28159.8542017977
$ws.Range("B5").Value = 484.3849148771696
$ws.Range("C5").Value = -2721.999989754528
$ws.Range("D5").Value = -2787.480470294639
$ws.Range("E5").Value = -2839.230548645945
$ws.Range("F5").Value = -2878.053769894919
$ws.Range("G5").Value = -3029.325695304055
$ws.Range("H5").Value = -3363.241479744785
$ws.Range("I5").Value = -3696.600133959575
$ws.Range("J5").Value = -5013.252792621915
$ws.Range("K5").Value = -5034.66390778905
$ws.Range("B6").Value = 1681.227758686155
$ws.Range("C6").Value = 1081.420427045613
$ws.Range("D6").Value = 1069.294750926665
$ws.Range("E6").Value = 1059.723614166635
$ws.Range("F6").Value = 1052.562796005634
$ws.Range("G6").Value = 1024.515715959518
$ws.Range("H6").Value = 962.5784025424464
$ws.Range("I6").Value = 900.7873297381775
$ws.Range("J6").Value = 651.5398652914637
$ws.Range("K6").Value = 637.3431019090867
$ws.Range("B7").Value = 1567.5
$ws.Range("C7").Value = 1567.5
$ws.Range("D7").Value = 1567.5
$ws.Range("E7").Value = 1567.5
$ws.Range("F7").Value = 1567.5
$ws.Range("G7").Value = 1567.5
$ws.Range("H7").Value = 1567.5
$ws.Range("I7").Value = 1567.5
$ws.Range("J7").Value = 1567.5
$ws.Range("K7").Value = 1567.5
$ws.Range("B8").Value = 113.7277586861551
$ws.Range("C8").Value = -486.079572954387
$ws.Range("D8").Value = -498.2052490733349
$ws.Range("E8").Value = -507.7763858333647
$ws.Range("F8").Value = -514.9372039943664
$ws.Range("G8").Value = -542.9842840404817
$ws.Range("H8").Value = -604.9215974575536
$ws.Range("I8").Value = -666.7126702618225
$ws.Range("J8").Value = -915.9601347085363
$ws.Range("K8").Value = -930.1568980909133
$ws.Range("B9").Value = -5.451086425781227
$ws.Range("C9").Value = -1.165106201171852
$ws.Range("D9").Value = -1.423010253906227
$ws.Range("E9").Value = -1.710668945312477
$ws.Range("F9").Value = -2.026281738281227
$ws.Range("G9").Value = -2.214483642578102
$ws.Range("H9").Value = -2.187811279296852
$ws.Range("I9").Value = -2.171179199218727
$ws.Range("J9").Value = -1.440283203124977
$ws.Range("K9").Value = -2.182440185546852
$ws.Range("B10").Value = 23.84453740855605
$ws.Range("C10").Value = 23.98732364919283
$ws.Range("D10").Value = 23.98801190270962
$ws.Range("E10").Value = 23.99132867835647
$ws.Range("F10").Value = 23.99714815349495
$ws.Range("G10").Value = 23.97237388782048
$ws.Range("H10").Value = 23.89839661419558
$ws.Range("I10").Value = 23.82359664696867
$ws.Range("J10").Value = 23.791955280438
$ws.Range("K10").Value = 24.00619918249913
$ws.Range("B11").Value = 24.00000000004019
$ws.Range("C11").Value = 23.9999999998825
$ws.Range("D11").Value = 23.99999999999557
$ws.Range("E11").Value = 23.99999999999926
$ws.Range("F11").Value = 23.99999999999983
$ws.Range("G11").Value = 23.9999999999996
$ws.Range("H11").Value = 23.99999999999972
$ws.Range("I11").Value = 23.99999999999977
$ws.Range("J11").Value = 23.99999999999983
$ws.Range("K11").Value = 23.99999999995293
$ws.Range("B12").Value = 24.10240300547554
$ws.Range("C12").Value = 23.58161050234139
$ws.Range("D12").Value = 23.57192269146827
$ws.Range("E12").Value = 23.56418304479172
$ws.Range("F12").Value = 23.5582666181997
$ws.Range("G12").Value = 23.53678143226557
$ws.Range("H12").Value = 23.4899411607214
$ws.Range("I12").Value = 23.44316849016877
$ws.Range("J12").Value = 23.25738205723979
$ws.Range("K12").Value = 23.12202466584654
$ws.Range("B13").Value = 24.04445599188716
$ws.Range("C13").Value = 23.18810105603154
$ws.Range("D13").Value = 23.16961254466679
$ws.Range("E13").Value = 23.15742066529475
$ws.Range("F13").Value = 23.15117223160536
$ws.Range("G13").Value = 23.08527498221935
$ws.Range("H13").Value = 22.92276118539974
$ws.Range("I13").Value = 22.75949828507527
$ws.Range("J13").Value = 22.37520044139393
$ws.Range("K13").Value = 22.37694212082602
$ws.Range("B14").Value = 12.23839114264644
$ws.Range("C14").Value = 14.6305321978881
$ws.Range("D14").Value = 14.67646779989116
$ws.Range("E14").Value = 14.71459397409853
$ws.Range("F14").Value = 14.74535530375152
$ws.Range("G14").Value = 14.83469630730582
$ws.Range("H14").Value = 15.01952440727166
$ws.Range("I14").Value = 15.20352381637599
$ws.Range("J14").Value = 15.91260021927872
$ws.Range("K14").Value = 14.21063626146469
$ws.Range("B15").Value = 9.212165898480206
$ws.Range("C15").Value = 8.166351844817301
$ws.Range("D15").Value = 7.920414662332316
$ws.Range("E15").Value = 7.685613563694098
$ws.Range("F15").Value = 7.461712185378011
$ws.Range("G15").Value = 7.262864872021299
$ws.Range("H15").Value = 7.095582782537317
$ws.Range("I15").Value = 6.937872541435467
$ws.Range("J15").Value = 6.550943243122166
$ws.Range("K15").Value = 7.52642136078125
$ws.Range("B16").Value = 17.30230340426689
$ws.Range("C16").Value = 15.63374472765234
$ws.Range("D16").Value = 15.1658840762182
$ws.Range("E16").Value = 14.7214294281179
$ws.Range("F16").Value = 14.29963744841218
$ws.Range("G16").Value = 13.89975993094914
$ws.Range("H16").Value = 13.52104983533501
$ws.Range("I16").Value = 13.16276575722204
$ws.Range("J16").Value = 12.46305065221908
$ws.Range("K16").Value = 14.77964047872252
$ws.Range("B17").Value = 17.16620932604729
$ws.Range("C17").Value = 16.00107577486694
$ws.Range("D17").Value = 15.53026989260349
$ws.Range("E17").Value = 15.08117311599967
$ws.Range("F17").Value = 14.65332348047074
$ws.Range("G17").Value = 14.26106214382601
$ws.Range("H17").Value = 13.91081364662703
$ws.Range("I17").Value = 13.57955419191072
$ws.Range("J17").Value = 13.0008867317437
$ws.Range("K17").Value = 15.55068751431124
$ws.Range("B18").Value = 9.100095945018122
$ws.Range("C18").Value = 8.56714512595348
$ws.Range("D18").Value = 8.318705261110853
$ws.Range("E18").Value = 8.079591769278982
$ws.Range("F18").Value = 7.849855194065776
$ws.Range("G18").Value = 7.659702948137785
$ws.Range("H18").Value = 7.523670832374512
$ws.Range("I18").Value = 7.396193354346611
$ws.Range("J18").Value = 7.134897818327829
$ws.Range("K18").Value = 8.303131772822105
$ws.Range("B19").Value = 1.99999999997838
$ws.Range("C19").Value = 2.120000000011354
$ws.Range("D19").Value = 2.120000000000585
$ws.Range("E19").Value = 2.120000000000123
$ws.Range("F19").Value = 2.12000000000002
$ws.Range("G19").Value = 2.120000000000012
$ws.Range("H19").Value = 2.120000000000016
$ws.Range("I19").Value = 2.120000000000026
$ws.Range("J19").Value = 2.120000000000017
$ws.Range("K19").Value = 1.760000000010415
$ws.Range("B20").Value = 0.2296064874682811
$ws.Range("C20").Value = 0.2450747733348159
$ws.Range("D20").Value = 0.2448723215702357
$ws.Range("E20").Value = 0.2446808965583268
$ws.Range("F20").Value = 0.244500038636872
$ws.Range("G20").Value = 0.2443278300019649
$ws.Range("H20").Value = 0.2441631498797155
$ws.Range("I20").Value = 0.2440078004336775
$ws.Range("J20").Value = 0.2446989847828288
$ws.Range("K20").Value = 0.2071529920507962
$ws.Range("B21").Value = 1.783393512510098
$ws.Range("C21").Value = 1.890925226676538
$ws.Range("D21").Value = 1.89112767843035
$ws.Range("E21").Value = 1.891319103441796
$ws.Range("F21").Value = 1.891499961363148
$ws.Range("G21").Value = 1.891672169998047
$ws.Range("H21").Value = 1.8918368501203
$ws.Range("I21").Value = 1.891992199566349
$ws.Range("J21").Value = 1.892301015217188
$ws.Range("K21").Value = 1.570847007959618
$ws.Range("B22").Value = 1.786359045506833
$ws.Range("C22").Value = 1.893918641626515
$ws.Range("D22").Value = 1.894128911433072
$ws.Range("E22").Value = 1.894327763373149
$ws.Range("F22").Value = 1.894515669525717
$ws.Range("G22").Value = 1.894694560196804
$ws.Range("H22").Value = 1.894865568642563
$ws.Range("I22").Value = 1.895026905089326
$ws.Range("J22").Value = 1.895347413124109
$ws.Range("K22").Value = 1.57385469517074
$ws.Range("B23").Value = 0.2296654433012009
$ws.Range("C23").Value = 0.2451339960098267
$ws.Range("D23").Value = 0.2449316084384918
$ws.Range("E23").Value = 0.2447402477264404
$ws.Range("F23").Value = 0.2445594519376755
$ws.Range("G23").Value = 0.244387298822403
$ws.Range("H23").Value = 0.2442226856946945
$ws.Range("I23").Value = 0.2440673857927322
$ws.Range("J23").Value = 0.2447586804628372
$ws.Range("K23").Value = 0.2072123885154724
$ws.Range("B25").Value = 1.016024470329285
$ws.Range("C25").Value = 1.079052686691284
$ws.Range("D25").Value = 1.079060554504395
$ws.Range("E25").Value = 1.079068064689636
$ws.Range("F25").Value = 1.07907509803772
$ws.Range("G25").Value = 1.079081892967224
$ws.Range("H25").Value = 1.07908821105957
$ws.Range("I25").Value = 1.079094290733337
$ws.Range("J25").Value = 1.080106139183044
$ws.Range("K25").Value = 0.9010670781135559
$ws.Range("B26").Value = 0.1214574831347294
$ws.Range("C26").Value = 0.1211442790561498
$ws.Range("D26").Value = 0.1210242566806496
$ws.Range("E26").Value = 0.1209107950859146
$ws.Range("F26").Value = 0.1208036179245803
$ws.Range("G26").Value = 0.1207015854138197
$ws.Range("H26").Value = 0.1206040308735961
$ws.Range("I26").Value = 0.1205120192810191
$ws.Range("J26").Value = 0.120329156382498
$ws.Range("K26").Value = 0.1204146496077222
$ws.Range("B29").Value = 104.5764609375
$ws.Range("C29").Value = 99.853296875
$ws.Range("D29").Value = 96.79503125
$ws.Range("E29").Value = 93.9093125
$ws.Range("F29").Value = 91.188640625
$ws.Range("G29").Value = 88.626421875
$ws.Range("H29").Value = 86.216046875
$ws.Range("I29").Value = 83.94959375
$ws.Range("J29").Value = 79.5993359375
$ws.Range("K29").Value = 78.571546875
$ws.Range("B30").Value = 100.1463671875
$ws.Range("C30").Value = 94.8722578125
$ws.Range("D30").Value = 91.7260625
$ws.Range("E30").Value = 88.7519453125
$ws.Range("F30").Value = 85.942609375
$ws.Range("G30").Value = 83.290984375
$ws.Range("H30").Value = 80.7901875
$ws.Range("I30").Value = 78.4335859375
$ws.Range("J30").Value = 73.89046875
$ws.Range("K30").Value = 74.2826171875
$ws.Range("B31").Value = 104.5764627205822
$ws.Range("C31").Value = 99.85329808712335
$ws.Range("D31").Value = 96.79502769967377
$ws.Range("E31").Value = 93.90931025889061
$ws.Range("F31").Value = 91.18863749321415
$ws.Range("G31").Value = 88.62642336877202
$ws.Range("H31").Value = 86.2160436511466
$ws.Range("I31").Value = 83.949590772445
$ws.Range("J31").Value = 79.59933254962218
$ws.Range("K31").Value = 78.57154913463265
